$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '50.893.42'
Set-TextValue 'E2' '  -1.98%  '

Set-TextValue 'D3' '2.932.41'
Set-TextValue 'E3' '  -3.13%  '

Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.16%  '

Set-TextValue 'D5' '375.87'
Set-TextValue 'E5' '  -1.82%  '

Set-TextValue 'D6' '101.88'
Set-TextValue 'E6' '  -4.71%  '

Set-TextValue 'E7' '  -3.04%  '

Set-TextValue 'D9' '0.581'
Set-TextValue 'E9' '  -4.11%  '

Set-TextValue 'D10' '36.42'
Set-TextValue 'E10' '  -4.51%  '

Set-TextValue 'D11' '0.139'
Set-TextValue 'E11' '  -1.18%  '

Set-TextValue 'D12' '0.0832'
Set-TextValue 'E12' '  -2.23%  '

Set-TextValue 'D13' '3.396.06'
Set-TextValue 'E13' '  -2.61%  '

Set-TextValue 'D14' '17.90'
Set-TextValue 'E14' '  -5.80%  '

Set-TextValue 'D15' '7.31'
Set-TextValue 'E15' '  -3.83%  '

Set-TextValue 'D16' '2.905.71'
Set-TextValue 'E16' '  -3.48%  '

Set-TextValue 'D17' '0.969'
Set-TextValue 'E17' '  -1.40%  '

Set-TextValue 'D18' '50.895.90'
Set-TextValue 'E18' '  -1.90%  '

Set-TextValue 'E19' '  -8.11%  '

Set-TextValue 'D20' '7.11'
Set-TextValue 'E20' '  -5.31%  '

Set-TextValue 'D21' '12.43'
Set-TextValue 'E21' '  -6.10%  '

Set-TextValue 'D22' '0.0₃0944'
Set-TextValue 'E22' '  -2.51%  '

Set-TextValue 'D23' '68.00'
Set-TextValue 'E23' '  -1.72%  '

Set-TextValue 'D24' '261.27'
Set-TextValue 'E24' '  -1.52%  '

Set-TextValue 'D25' '2.85'
Set-TextValue 'E25' '  +1.02%  '

Set-TextValue 'D26' '8.19'
Set-TextValue 'E26' '  +7.89%  '

Set-TextValue 'D27' '7.67'
Set-TextValue 'E27' '  +4.34%  '

Set-TextValue 'D28' '0.167'
Set-TextValue 'E28' '  -4.07%  '

Set-TextValue 'E29' '  +0.00%  '

Set-TextValue 'E30' '  +4.75%  '

Set-TextValue 'D31' '25.55'
Set-TextValue 'E31' '  -3.20%  '

Set-TextValue 'D32' '9.78'
Set-TextValue 'E32' '  -2.33%  '

Set-TextValue 'D33' '0.0457'
Set-TextValue 'E33' '  +1.34%  '

Set-TextValue 'D34' '50.59'
Set-TextValue 'E34' '  -1.22%  '

Set-TextValue 'D35' '33.76'
Set-TextValue 'E35' '  -4.23%  '

Set-TextValue 'E36' '  -3.16%  '

Set-TextValue 'E37' '  +0.19%  '

Set-TextValue 'D38' '2.96'
Set-TextValue 'E38' '  -5.79%  '

Set-TextValue 'D39' '2.54'
Set-TextValue 'E39' '  -5.04%  '

Set-TextValue 'E40' '  -2.82%  '

Set-TextValue 'D41' '16.27'
Set-TextValue 'E41' '  -8.29%  '

Set-TextValue 'D42' '1.77'
Set-TextValue 'E42' '  -6.15%  '

Set-TextValue 'D43' '121.08'
Set-TextValue 'E43' '  -3.29%  '

Set-TextValue 'D44' '21.01'
Set-TextValue 'E44' '  -7.09%  '

Set-TextValue 'D45' '2.04'
Set-TextValue 'E45' '  -2.17%  '

Set-TextValue 'E46' '  -1.72%  '

Set-TextValue 'D47' '0.270'
Set-TextValue 'E47' '  -3.10%  '

Set-TextValue 'D48' '2.001.10'
Set-TextValue 'E48' '  -3.23%  '

Set-TextValue 'D49' '3.20'
Set-TextValue 'E49' '  -3.78%  '

Set-TextValue 'D50' '0.0340'
Set-TextValue 'E50' '  -4.40%  '

Set-TextValue 'D51' '4.99'
Set-TextValue 'E51' '  -5.12%  '

